$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.954.12"
$ws.Range("E2").Value = "  +5.17%  "

$ws.Range("D3").Value = "3.247.15"
$ws.Range("E3").Value = "  +2.67%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'395.11"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").Value = "'108.16"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("D7").Value = "'0.586"
$ws.Range("E7").Value = "  +6.94%  "

$ws.Range("D8").Value = "3.241.95"
$ws.Range("E8").Value = "  +2.73%  "

$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").Value = "'0.625"
$ws.Range("E10").Value = "  +1.53%  "

$ws.Range("D11").Value = "'39.17"
$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").Value = "'0.0977"
$ws.Range("E12").Value = "  +11.21%  "

$ws.Range("E13").Value = "  +1.65%  "

$ws.Range("D14").Value = "3.755.80"
$ws.Range("E14").Value = "  +3.39%  "

$ws.Range("D15").Value = "'8.21"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").Value = "'19.14"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "3.217.38"
$ws.Range("E17").Value = "  +2.41%  "

$ws.Range("E18").Value = "  -2.46%  "

$ws.Range("D19").Value = "'10.85"
$ws.Range("E19").Value = "  +3.54%  "

$ws.Range("D20").Value = "56.847.82"
$ws.Range("E20").Value = "  +5.76%  "

$ws.Range("D21").Value = "'3.36"
$ws.Range("E21").Value = "  +2.58%  "

$ws.Range("E22").Value = "  +12.32%  "

$ws.Range("E23").Value = "  +1.58%  "

$ws.Range("D24").Value = "'296.17"
$ws.Range("E24").Value = "  +9.09%  "

$ws.Range("D25").Value = "'74.27"
$ws.Range("E25").Value = "  +4.57%  "

$ws.Range("E26").Value = "  -3.02%  "

$ws.Range("D27").Value = "'27.83"
$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("D28").Value = "'7.66"
$ws.Range("E28").Value = "  -3.99%  "

$ws.Range("D29").Value = "'7.34"
$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("E32").Value = "  +4.11%  "

$ws.Range("E33").Value = "  -2.18%  "

$ws.Range("D34").Value = "'39.35"
$ws.Range("E34").Value = "  +6.45%  "

$ws.Range("D35").Value = "'0.0482"
$ws.Range("E35").Value = "  -4.25%  "

$ws.Range("E36").Value = "  +1.70%  "

$ws.Range("D37").Value = "'51.75"
$ws.Range("E37").Value = "  +2.55%  "

$ws.Range("D38").Value = "'3.51"
$ws.Range("E38").Value = "  -3.75%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("E40").Value = "  +3.28%  "

$ws.Range("D41").Value = "'135.34"
$ws.Range("E41").Value = "  +3.84%  "

$ws.Range("E42").Value = "  +3.79%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'17.09"
$ws.Range("E43").Value = "  -1.28%  "

$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.94"
$ws.Range("E44").Value = "  -3.81%  "

$ws.Range("D45").Value = "'1.89"
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("E46").Value = "  -2.60%  "

$ws.Range("D47").Value = "'22.19"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").Value = "'2.13"
$ws.Range("E48").Value = "  +2.31%  "

$ws.Range("D49").Value = "2.162.93"
$ws.Range("E49").Value = "  +4.05%  "

$ws.Range("D50").Value = "'1.98"
$ws.Range("E50").Value = "  +18.53%  "

$ws.Range("E51").Value = "  -3.50%  "
